# Applies the Thu Mar 28 13:06:36 UTC 2024 cryptos-list refresh (GitHub Actions job).
# Updates price (col D) and 1h-volume (col E) figures for each coin row, plus the
# two rank swaps (Toncoin/PancakeSwap at rows 24-25, Hedera/OKB at rows 33-34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.800.67"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.576.24"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.33"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.57"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.568.94"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.216"
$ws.Range("E10").Value = "  +14.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.649"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.24"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("E13").Value = "  +5.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.51"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.140.38"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.37"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.709.05"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.558.35"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "578.01"
$ws.Range("E19").Value = "  +16.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.38"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.01"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.76"
$ws.Range("E23").Value = "  -8.86%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.04"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.59"
$ws.Range("E25").Value = "  +4.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.34"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.93"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.12"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.23"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.33"
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.30"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.115"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.18"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.30"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "555.42"
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.412"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0813"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.58"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.429.27"
$ws.Range("E41").Value = "  +6.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.137"
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.10"
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.49"
$ws.Range("E45").Value = "  -6.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0444"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.96"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.39"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.138"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.997"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  -5.71%  "
